$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite rows 2-3 with updated (rounded) model output and add three
# brand-new rows (4-6) for the newly enabled models -- 5 users total now.
$data = @(
    @(1, 0.56, 0.59, 0, -9.6, 0.09, 0.329, 0.18, 0.182, 0.42, 121, 4, 5, 8, 8, 5, 5),
    @(2, 0.61, 0.67, 0, -6.6, 0.1, 0.21, 0.04, 0.165, 0.5, 125, 4, 8, 6, 6, 9, 10),
    @(3, 0.57, 0.75, 1, -5.9, 0.09, 0.173, 0.06, 0.198, 0.49, 127, 4, 3, 8, 4, 8, 8),
    @(4, 0.56, 0.54, 0, -10.7, 0.09, 0.376, 0.23, 0.173, 0.39, 119, 4, 5, 8, 8, 5, 5),
    @(5, 0.83, 0.65, 6, -5.2, 0.06, 0.563, 0, 0.097, 0.34, 136, 4, 6, 7, 7, 6, 6)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

# Reset page margins to Excel defaults (in points; Excel COM uses points)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Reset view: select B7 (also clears the stale topLeftCell scroll state)
[void]$ws.Range("B7").Select()
